# Fix class name on the worksheet: "CE En 544" -> "CE 544"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A2").Value = "CE 544 - Brigham Young University"
